# Appends a new row (row 78) to each of the four worksheets, mirroring the
# existing row 77 on that sheet but with the timestamp advanced by one hour.

$wb = $excel.ActiveWorkbook

$rows = @(
    @{ Sheet = "ROW35-FE-LIFTER";  A = "2025-03-07 13:42:06"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"; E = "0x d";  G = "568631262647113770877196"; I = 13 },
    @{ Sheet = "ROW35-MID-LIFTER"; A = "2025-03-07 13:29:35"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"; E = "0x e";  G = "568631262647113770942732"; I = 14 },
    @{ Sheet = "ROW02-FE-LIFTER";  A = "2025-03-07 13:51:45"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x06,0x41,0x0c,"; E = "0xff";  G = "568631262647113769959692"; I = 255 },
    @{ Sheet = "ROW02-MID-LIFTER"; A = "2025-03-07 13:41:15"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"; E = "0x 3";  G = "568631262647113769959692"; I = 3 }
)

foreach ($rowInfo in $rows) {
    $ws = $wb.Worksheets.Item($rowInfo.Sheet)
    $newRow = 78

    $ws.Cells.Item($newRow, 1).Value = $rowInfo.A
    $ws.Cells.Item($newRow, 2).Value = "0x01,0x90 "
    $ws.Cells.Item($newRow, 3).Value = $rowInfo.C
    $ws.Cells.Item($newRow, 4).Value = "0x01,0x90,"
    $ws.Cells.Item($newRow, 5).Value = $rowInfo.E
    $ws.Cells.Item($newRow, 6).Value = 400

    # Column G is a large integer that exceeds double precision, so it must
    # stay stored as text (matching every other row on these sheets).
    $ws.Cells.Item($newRow, 7).NumberFormat = "@"
    $ws.Cells.Item($newRow, 7).Value = $rowInfo.G

    $ws.Cells.Item($newRow, 8).Value = 400
    $ws.Cells.Item($newRow, 9).Value = $rowInfo.I
}
